$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> row -> new value for column F ("想去人数")
$sheetNames = @("展览", "全部类型")

$updates = @{
    6  = 521
    8  = 1970
    10 = 87
    11 = 4160
    15 = 87
    18 = 2885
    20 = 405
    25 = 59
    30 = 305
    31 = 1634
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
